$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -le 19; $i++) {
    $row = $i + 2
    $suffix = "{0:D2}" -f $i
    $ws.Range("A$row").Value = "sequences/278857_learning_sequence_$suffix.csv"
}
